# "New Sheets Alert!" update — merge the pivot's "New" (C) and "Present" (D)
# columns on the "Status by State" sheet into a single "Present" column,
# then drop the now-redundant D column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Status by State")

# Last used row in column A (header is row 1, data starts row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $newVal = $ws.Cells.Item($r, 3).Value2
    $presentVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value = $newVal + $presentVal
}

# Column C's header becomes "Present" (absorbing column D's former header).
$ws.Range("C1").Value = "Present"

# Column D ("Present") is now folded into C — delete it, shifting nothing else.
$ws.Columns("D").Delete()
